# DuplicatedCharsSimpleValues.docx edit:
#   "{abcdef}"            ->  "{xyz}"
#   "third {a{u}, ..."    ->  "third {a{xyz}, ..."
#   "fourth {f}d}"        ->  "fourth {xyz}d}"
# Each replaced word ("abcdef" / "u" / "f") must land in its own <w:r>,
# with the surrounding literal text kept split into the same runs Word
# already used for it, exactly as the diff shows.
#
# This engine normalizes (merges) any two adjacent runs that end up
# with identical rPr whenever an edit touches the paragraph, which
# would silently re-glue the paragraph's existing run splits back
# together the moment we change any text in it. To stop that, every
# segment of the paragraph that must remain its own run is first
# "fenced" with a tiny, reversible formatting nudge (Bold, alternating
# on/off between neighbours) so adjacent segments never share identical
# rPr while we are still cutting things up; the nudge is removed again
# (Bold = 0 fully deletes the <w:b/> element, unlike e.g. Italic which
# leaves a val="0" residue) only once every split has been made and
# every replacement typed in.

$d = $word.ActiveDocument
$full = $d.Content.Text.TrimEnd([char]13, [char]7)

# --- split the paragraph's current text into its existing run pieces --
# (mirrors the <w:t> runs already in word/document.xml for this
# paragraph) plus the three words that are about to be replaced.
$literalSegments = @(
    "Firs", "t", ": {", "{", "xyz}", ", second", ": ", "{",
    "abcdef", "}", "}, ", "third ", "{", "a{", "u", "}",
    ", fourth {", "f", "}d}"
)

# --- compute [start,end) for every segment by walking the text --------
$bounds = New-Object 'System.Object[][]' $literalSegments.Length
$cursor = 0
for ($i = 0; $i -lt $literalSegments.Length; $i++) {
    $seg = $literalSegments[$i]
    $start = $full.IndexOf($seg, $cursor)
    if ($start -ne $cursor) { throw "segment '$seg' did not start at cursor ($cursor vs $start)" }
    $end = $start + $seg.Length
    $bounds[$i] = @($start, $end)
    $cursor = $end
}
if ($cursor -ne $full.Length) { throw "segments did not cover the whole paragraph ($cursor vs $($full.Length))" }

# --- fence every boundary, alternating so neighbours always differ ----
for ($i = 0; $i -lt $bounds.Length; $i++) {
    if ($i % 2 -eq 1) {
        $s = $bounds[$i]
        $d.Range($s[0], $s[1]).Bold = 1
    }
}

# --- replace the three words, right-to-left so earlier offsets stay --
# valid (indices 17, 14, 8 -- "f", "u", "abcdef" in that order)
$s17 = $bounds[17]; $d.Range($s17[0], $s17[1]).Text = "xyz"
$s14 = $bounds[14]; $d.Range($s14[0], $s14[1]).Text = "xyz"
$s8 = $bounds[8];   $d.Range($s8[0], $s8[1]).Text = "xyz"

# --- drop the fencing; Bold=0 removes <w:b/> outright (no val="0") ----
$d.Range(0, $d.Content.End).Bold = 0

$d.Content.Text
